# daily auto push: 2026-02-13 09:58 UTC
# Insert a new data row for 2026/02/13 (Fri) above the existing row 814,
# shifting all subsequent rows (old 814:855) down to (815:856).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(814).Insert()

# Use a leading apostrophe so Excel stores the date-looking value as plain
# text (matching how every other row in column A is stored), then reset
# the style back to Normal so no stray formatting/style index lingers.
$ws.Cells.Item(814, 1).Value = "'2026/02/13"
$ws.Cells.Item(814, 1).Style = "Normal"
$ws.Cells.Item(814, 2).Value = "金"
$ws.Cells.Item(814, 3).Value = 16
$ws.Cells.Item(814, 4).Value = 201
